$d = $word.ActiveDocument

# 1. Change the first paragraph's style from Heading1 to Title
$d.Paragraphs(1).Style = "Title"

# 2. Update the built-in "Title" style definition:
#    spacing-after goes from 60 twips (3pt) to 240 twips (12pt)
$titleStyle = $d.Styles("Title")
$titleStyle.ParagraphFormat.SpaceAfter = 12
